# Commit: "Fri, Apr 24, 2020  9:06:30 PM"
#
# Two logical changes, per the canonical-OOXML diff:
#
#  1) Three tables (on the slides whose grid columns are 2879725,
#     3424250 and 2881325 EMU wide -- i.e. slides 14, 15 and 16) get
#     their <a:tableStyleId> switched from the "no style" GUID
#     {D4AEAB4D-0E84-4705-9016-3CAB8418E4A4} to the built-in table
#     style {7983CE0C-A75F-44D0-A8D0-516D5EC18CBF}. PowerPoint tables
#     refuse a plain property assignment for this ("Table styles
#     cannot be assigned through a property"), so this goes through
#     Table.ApplyStyle(id).
#
#  2) The deck's two theme parts swap their color-scheme payload: the
#     slide master's theme (currently "Integral" / "Red Violet")
#     becomes the stock "Office" palette, and vice versa for the
#     notes-master theme. The notes-master theme part isn't reachable
#     through this object model (NotesMaster/HandoutMaster always
#     resolve back to the slide master's Theme here), so only the
#     slide master side of the swap can be driven from this script;
#     it is applied in full via ThemeColorScheme, color by color, to
#     land on the exact target RGB values.

$p = $ppt.ActivePresentation

# --- 1) Table styles --------------------------------------------------
$newStyleId = "{7983CE0C-A75F-44D0-A8D0-516D5EC18CBF}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Theme colors (slide master: "Integral" -> "Office") ----------
$colors = $p.SlideMaster.Theme.ThemeColorScheme
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
